$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 29.75
$ws.Range("I11").Value = 29.75
$ws.Range("K11").Value = 29.75
$ws.Range("M11").Value = 110.25
$ws.Range("H18").Value = 835.7143
$ws.Range("I18").Value = 725
$ws.Range("K18").Value = 725
$ws.Range("M18").Value = -441
$ws.Range("H53").Value = 488.6154
$ws.Range("I53").Value = 496.2
$ws.Range("K53").Value = 496.2
$ws.Range("M53").Value = 140.8
$ws.Range("H58").Value = 3186
$ws.Range("J58").Value = 4231.3335
$ws.Range("L58").Value = 12694.0005
$ws.Range("N58").Value = -12994.0005
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H98").Value = 575
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 1000
$ws.Range("N98").Value = -3996
$ws.Range("H122").Value = 575
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 1459.6
$ws.Range("I132").Value = 1459.6
$ws.Range("K132").Value = 4378.799999999999
$ws.Range("M132").Value = -1848.799999999999
$ws.Range("H135").Value = 1688.8334
$ws.Range("I135").Value = 1899.5
$ws.Range("J135").Value = 1267.5
$ws.Range("K135").Value = 17095.5
$ws.Range("L135").Value = 11407.5
$ws.Range("M135").Value = -14560.5
$ws.Range("N135").Value = -16477.5
$ws.Range("H138").Value = 2379.3333
$ws.Range("J138").Value = 2625
$ws.Range("L138").Value = 7875
$ws.Range("N138").Value = -18155

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 225
$ws.Range("I110").Value = 240
$ws.Range("J110").Value = 180
$ws.Range("K110").Value = 240
$ws.Range("L110").Value = 180
$ws.Range("M110").Value = 1805
$ws.Range("N110").Value = -4270

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H88").Value = 15699.8
$ws.Range("J88").Value = 15699.8
$ws.Range("L88").Value = 15699.8
$ws.Range("N88").Value = -16511.8
$ws.Range("H91").Value = 15699.8
$ws.Range("J91").Value = 15699.8
$ws.Range("L91").Value = 15699.8
$ws.Range("N91").Value = -18507.8
$ws.Range("H99").Value = 556
$ws.Range("I99").Value = 445
$ws.Range("K99").Value = 445
$ws.Range("M99").Value = 1053
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H10").Value = 3024.2856
$ws.Range("I10").Value = 193.66667
$ws.Range("J10").Value = 20008
$ws.Range("K10").Value = 193.66667
$ws.Range("L10").Value = 20008
$ws.Range("M10").Value = -54.66667000000001
$ws.Range("N10").Value = -20286
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700
$ws.Range("H35").Value = 233.57143
$ws.Range("I35").Value = 222.5
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 222.5
$ws.Range("L35").Value = 300
$ws.Range("M35").Value = 71.5
$ws.Range("N35").Value = -888
$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 7000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5877
$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 7000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 35000
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -29384
$ws.Range("H110").Value = 75000
$ws.Range("J110").Value = 75000
$ws.Range("L110").Value = 75000
$ws.Range("N110").Value = -83180

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H39").Value = 7921.4287
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 7921.4287
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 23764.2861
$ws.Range("N39").Value = -24352.2861
$ws.Range("H41").Value = 300
$ws.Range("I41").Value = 300
$ws.Range("K41").Value = 900
$ws.Range("M41").Value = -562
$ws.Range("H48").Value = 2368.3333
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H63").Value = 1500
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 4500
$ws.Range("M63").Value = -3751
$ws.Range("H66").Value = 1500
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 13500
$ws.Range("M66").Value = -9756
$ws.Range("H69").Value = 1551
$ws.Range("I69").Value = 1012
$ws.Range("J69").Value = 1658.8
$ws.Range("K69").Value = 3036
$ws.Range("L69").Value = 4976.4
$ws.Range("M69").Value = -2225
$ws.Range("N69").Value = -6598.4
$ws.Range("H72").Value = 1551
$ws.Range("I72").Value = 1012
$ws.Range("J72").Value = 1658.8
$ws.Range("K72").Value = 9108
$ws.Range("L72").Value = 14929.2
$ws.Range("M72").Value = -5052
$ws.Range("N72").Value = -23041.2
$ws.Range("H131").Value = 1580.2
$ws.Range("I131").Value = 892.5
$ws.Range("J131").Value = 1752.125
$ws.Range("K131").Value = 2677.5
$ws.Range("L131").Value = 5256.375
$ws.Range("M131").Value = 2362.5
$ws.Range("N131").Value = -15336.375
$ws.Range("H132").Value = 2447.75
$ws.Range("I132").Value = 1897
$ws.Range("J132").Value = 2998.5
$ws.Range("K132").Value = 17073
$ws.Range("L132").Value = 26986.5
$ws.Range("M132").Value = -14543
$ws.Range("N132").Value = -32046.5

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 5301.353
$ws.Range("I113").Value = 5010.5713
$ws.Range("J113").Value = 6658.3335
$ws.Range("K113").Value = 5010.5713
$ws.Range("L113").Value = 6658.3335
$ws.Range("M113").Value = -2840.5713
$ws.Range("N113").Value = -10998.3335

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -5
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 300
$ws.Range("K27").Value = 300
$ws.Range("M27").Value = -193
$ws.Range("H40").Value = 6033.1665
$ws.Range("I40").Value = 3460.6
$ws.Range("K40").Value = 3460.6
$ws.Range("M40").Value = -3324.6
$ws.Range("H111").Value = 71387
$ws.Range("J111").Value = 71387
$ws.Range("L111").Value = 71387
$ws.Range("N111").Value = -79567
$ws.Range("H122").Value = 2659.6667
$ws.Range("J122").Value = 3650
$ws.Range("L122").Value = 10950
$ws.Range("N122").Value = -15850

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 2282.3333
$ws.Range("I122").Value = 2282.3333
$ws.Range("K122").Value = 6846.999899999999
$ws.Range("M122").Value = -4396.999899999999
